# Updates cryptos list prices/volume deltas (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, even when it looks like a number
# (e.g. "1.00", "0.999"), without leaving the cell number-formatted -
# matches the original inlineStr/shared-string "@"-free cells.
function Set-TextValue($range, [string]$text) {
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Looks like a genuine number to Excel's parser - force text entry
        # via a leading apostrophe, then strip the resulting quote-prefix
        # style back to Normal so no extra formatting is left behind.
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value2 = $text
    }
}

Set-TextValue $ws.Range("D2") '37.008.09'
Set-TextValue $ws.Range("E2") '  -1.59%  '
Set-TextValue $ws.Range("D3") '2.014.91'
Set-TextValue $ws.Range("E3") '  -2.84%  '
Set-TextValue $ws.Range("E4") '  -0.04%  '
Set-TextValue $ws.Range("D5") '226.04'
Set-TextValue $ws.Range("E5") '  -2.49%  '
Set-TextValue $ws.Range("E6") '  -2.56%  '
Set-TextValue $ws.Range("E7") '  +0.08%  '
Set-TextValue $ws.Range("D8") '54.86'
Set-TextValue $ws.Range("E8") '  -5.46%  '
Set-TextValue $ws.Range("D9") '0.376'
Set-TextValue $ws.Range("E9") '  -2.94%  '
Set-TextValue $ws.Range("D10") '0.0781'
Set-TextValue $ws.Range("E10") '  +0.83%  '
Set-TextValue $ws.Range("D12") '2.313.84'
Set-TextValue $ws.Range("E12") '  -2.64%  '
Set-TextValue $ws.Range("D13") '14.12'
Set-TextValue $ws.Range("E13") '  -4.35%  '
Set-TextValue $ws.Range("D14") '20.21'
Set-TextValue $ws.Range("E14") '  -4.68%  '
Set-TextValue $ws.Range("E15") '  -3.14%  '
Set-TextValue $ws.Range("E16") '  -3.75%  '
Set-TextValue $ws.Range("D17") '2.039.63'
Set-TextValue $ws.Range("E17") '  -1.80%  '
Set-TextValue $ws.Range("D18") '37.019.51'
Set-TextValue $ws.Range("E18") '  -1.38%  '
Set-TextValue $ws.Range("D19") '6.18'
Set-TextValue $ws.Range("E19") '  +0.54%  '
Set-TextValue $ws.Range("D20") '68.81'
Set-TextValue $ws.Range("E20") '  -1.73%  '
Set-TextValue $ws.Range("E21") '  -1.38%  '
Set-TextValue $ws.Range("D22") '222.62'
Set-TextValue $ws.Range("E22") '  -2.02%  '
Set-TextValue $ws.Range("D23") '0.999'
Set-TextValue $ws.Range("E23") '  -0.09%  '
Set-TextValue $ws.Range("E24") '  +1.50%  '
Set-TextValue $ws.Range("E25") '  -6.49%  '
Set-TextValue $ws.Range("D26") '166.30'
Set-TextValue $ws.Range("E26") '  -1.94%  '
Set-TextValue $ws.Range("D27") '9.15'
Set-TextValue $ws.Range("E27") '  -7.39%  '
Set-TextValue $ws.Range("E28") '  -0.52%  '
Set-TextValue $ws.Range("B29") 'Kaspa'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D29") '0.124'
Set-TextValue $ws.Range("E29") '  -4.76%  '
Set-TextValue $ws.Range("B30") 'EthereumClassic'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D30") '18.66'
Set-TextValue $ws.Range("E30") '  -3.69%  '
Set-TextValue $ws.Range("E31") '  -3.83%  '
Set-TextValue $ws.Range("E32") '  -1.05%  '
Set-TextValue $ws.Range("D33") '0.0613'
Set-TextValue $ws.Range("E33") '  -2.19%  '
Set-TextValue $ws.Range("E34") '  -4.77%  '
Set-TextValue $ws.Range("E35") '  -7.22%  '
Set-TextValue $ws.Range("E36") '  +0.42%  '
Set-TextValue $ws.Range("D37") '1.00'
Set-TextValue $ws.Range("E37") '  +0.10%  '
Set-TextValue $ws.Range("E38") '  -4.41%  '
Set-TextValue $ws.Range("D39") '5.28'
Set-TextValue $ws.Range("E39") '  -0.66%  '
Set-TextValue $ws.Range("D40") '1.479.99'
Set-TextValue $ws.Range("E40") '  -0.45%  '
Set-TextValue $ws.Range("D41") '0.0215'
Set-TextValue $ws.Range("E41") '  -5.13%  '
Set-TextValue $ws.Range("E42") '  -3.39%  '
Set-TextValue $ws.Range("D43") '0.0915'
Set-TextValue $ws.Range("E43") '  -4.21%  '
Set-TextValue $ws.Range("D44") '16.30'
Set-TextValue $ws.Range("E44") '  -1.87%  '
Set-TextValue $ws.Range("E45") '  -5.08%  '
Set-TextValue $ws.Range("D46") '1.12'
Set-TextValue $ws.Range("E46") '  -5.43%  '
Set-TextValue $ws.Range("B47") 'ARBITRUM'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D47") '1.01'
Set-TextValue $ws.Range("E47") '  -2.77%  '
Set-TextValue $ws.Range("B48") 'FraxShare'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D48") '7.14'
Set-TextValue $ws.Range("E48") '  -1.62%  '
Set-TextValue $ws.Range("E49") '  -1.39%  '
Set-TextValue $ws.Range("D50") '2.201.51'
Set-TextValue $ws.Range("E50") '  -2.70%  '
Set-TextValue $ws.Range("D51") '44.37'
Set-TextValue $ws.Range("E51") '  -2.98%  '
